$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 1: date serial bumped by one day (2024-01-17 -> 2024-01-18)
$ws.Range("A1").Value = 45309

# Price updates (rows 32-37, column D)
$ws.Range("D32").Value = 7320
$ws.Range("D33").Value = 8170
$ws.Range("D34").Value = 9280
$ws.Range("D35").Value = 11550
$ws.Range("D36").Value = 11960
$ws.Range("D37").Value = 12590
